$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"999.62830350990407"
$ws.Range("D2").Value = [double]"7.4526371185434925e-11"
$ws.Range("B3").Value = [double]"16287.765014648438"
$ws.Range("D3").Value = [double]"9.8931515757350041e-11"
$ws.Range("B4").Value = [double]"115276.1064453125"
$ws.Range("D4").Value = [double]"5.8170468442142464e-10"
$ws.Range("B5").Value = [double]"572392.23046875"
$ws.Range("D5").Value = [double]"2.5505166867390017e-09"
$ws.Range("B6").Value = [double]"1459755.578125"
$ws.Range("D6").Value = [double]"1.661871529279324e-08"
$ws.Range("B7").Value = [double]"2317277.375"
$ws.Range("D7").Value = [double]"2.7449248207744859e-08"
$ws.Range("B8").Value = [double]"3562152.75"
$ws.Range("D8").Value = [double]"5.7005454578984427e-08"
$ws.Range("B9").Value = [double]"4595849.375"
$ws.Range("D9").Value = [double]"5.4672842253467024e-08"
$ws.Range("B10").Value = [double]"3743189.25"
$ws.Range("D10").Value = [double]"3.3145905575793222e-08"
$ws.Range("B11").Value = [double]"4709251.5"
$ws.Range("D11").Value = [double]"2.4011152177649819e-08"
$ws.Range("B12").Value = [double]"2928522.1875"
$ws.Range("D12").Value = [double]"3.0533104933283539e-08"
$ws.Range("B13").Value = [double]"2234498.9375"
$ws.Range("D13").Value = [double]"2.5565391581494623e-08"
$ws.Range("B14").Value = [double]"1100657.875"
$ws.Range("D14").Value = [double]"2.2078300077055246e-08"
$ws.Range("B15").Value = [double]"1142042.125"
$ws.Range("D15").Value = [double]"2.3925945669134308e-08"
$ws.Range("B16").Value = [double]"488321.96875"
$ws.Range("D16").Value = [double]"1.2484131950429855e-08"
$ws.Range("B17").Value = [double]"443130.1640625"
$ws.Range("D17").Value = [double]"2.1161328689345282e-08"
$ws.Range("B18").Value = [double]"164270.7578125"
$ws.Range("D18").Value = [double]"5.6076858712117428e-09"
$ws.Range("B19").Value = [double]"98320.80859375"
$ws.Range("D19").Value = [double]"4.7255777069210581e-09"
$ws.Range("B20").Value = [double]"56298.671875"
$ws.Range("D20").Value = [double]"2.6263191621467286e-09"
$ws.Range("B21").Value = [double]"30840.2724609375"
$ws.Range("D21").Value = [double]"3.1529672206431769e-09"
$ws.Range("B22").Value = [double]"16162.384765625"
$ws.Range("D22").Value = [double]"1.6369821054240674e-09"
$ws.Range("B23").Value = [double]"8103.25732421875"
$ws.Range("D23").Value = [double]"1.7871811808234384e-09"
$ws.Range("B24").Value = [double]"2307.703369140625"
$ws.Range("D24").Value = [double]"5.0268217366422618e-10"
$ws.Range("B25").Value = [double]"2647.4547119140625"
$ws.Range("D25").Value = [double]"1.6838178629186018e-09"
$ws.Range("B26").Value = [double]"1188.2264099121094"
$ws.Range("D26").Value = [double]"1.9252255345492131e-09"
$ws.Range("B27").Value = [double]"309.7191162109375"
$ws.Range("D27").Value = [double]"7.5273431932032508e-10"
